$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("B18").Value = "[-, -, 'MEC-2NA-Fundição', -]"
$ws.Range("D18").Value = "['MEC-2NA-Fundição', -, -, -]"
$ws.Range("E18").Value = "-"

# Row 19
$ws.Range("E19").Value = "-"

# Row 20
$ws.Range("E20").Value = "-"

# Row 21
$ws.Range("B21").Value = "[-, 'MEC-2NA-Fundição', -, -]"
$ws.Range("D21").Value = "['MEC-2NA-Fundição', -, -, -]"
$ws.Range("E21").Value = "-"
